$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1924.8837
$ws.Cells.Item(112, 9).Value = 828.5714
$ws.Cells.Item(112, 10).Value = 2138.0557
$ws.Cells.Item(112, 11).Value = 2485.7142
$ws.Cells.Item(112, 12).Value = 6414.1671
$ws.Cells.Item(112, 13).Value = -1377.7142
$ws.Cells.Item(112, 14).Value = -8630.167099999999
$ws.Cells.Item(129, 8).Value = 878.78845
$ws.Cells.Item(129, 9).Value = 319.85715
$ws.Cells.Item(129, 10).Value = 1084.7106
$ws.Cells.Item(129, 11).Value = 959.5714499999999
$ws.Cells.Item(129, 12).Value = 3254.1318
$ws.Cells.Item(129, 13).Value = 4040.42855
$ws.Cells.Item(129, 14).Value = -13254.1318
$ws.Cells.Item(141, 8).Value = 549.0526
$ws.Cells.Item(141, 9).Value = 549.0526
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 1647.1578
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = 3532.8422
$ws.Cells.Item(141, 14).ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4322.463
$ws.Cells.Item(32, 9).Value = 3991.2585
$ws.Cells.Item(32, 11).Value = 3991.2585
$ws.Cells.Item(32, 13).Value = -3704.2585
$ws.Cells.Item(45, 8).Value = 1351.2778
$ws.Cells.Item(45, 9).Value = 1323
$ws.Cells.Item(45, 11).Value = 1323
$ws.Cells.Item(45, 13).Value = -946
$ws.Cells.Item(135, 8).Value = 5427.8335
$ws.Cells.Item(135, 10).Value = 5427.8335
$ws.Cells.Item(135, 12).Value = 5427.8335
$ws.Cells.Item(135, 14).Value = -15567.8335
$ws.Cells.Item(139, 8).Value = 30032.5
$ws.Cells.Item(139, 10).Value = 30032.5
$ws.Cells.Item(139, 12).Value = 30032.5
$ws.Cells.Item(139, 14).Value = -40312.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3570.4666
$ws.Cells.Item(86, 9).Value = 4215.905
$ws.Cells.Item(86, 10).Value = 2064.4443
$ws.Cells.Item(86, 11).Value = 4215.905
$ws.Cells.Item(86, 12).Value = 2064.4443
$ws.Cells.Item(86, 13).Value = -3092.905
$ws.Cells.Item(86, 14).Value = -4310.4443
$ws.Cells.Item(89, 8).Value = 3570.4666
$ws.Cells.Item(89, 9).Value = 4215.905
$ws.Cells.Item(89, 10).Value = 2064.4443
$ws.Cells.Item(89, 11).Value = 21079.525
$ws.Cells.Item(89, 12).Value = 10322.2215
$ws.Cells.Item(89, 13).Value = -15463.525
$ws.Cells.Item(89, 14).Value = -21554.2215
$ws.Cells.Item(94, 8).Value = 22728654
$ws.Cells.Item(94, 9).Value = 25001218
$ws.Cells.Item(94, 10).Value = 3010
$ws.Cells.Item(94, 11).Value = 25001218
$ws.Cells.Item(94, 12).Value = 3010
$ws.Cells.Item(94, 13).Value = -25000767
$ws.Cells.Item(94, 14).Value = -3912
$ws.Cells.Item(134, 8).Value = 12327.583
$ws.Cells.Item(134, 9).Value = 1986.2
$ws.Cells.Item(134, 10).Value = 19714.285
$ws.Cells.Item(134, 11).Value = 5958.6
$ws.Cells.Item(134, 12).Value = 59142.855
$ws.Cells.Item(134, 13).Value = -3423.6
$ws.Cells.Item(134, 14).Value = -64212.855
$ws.Cells.Item(138, 8).Value = 50253.332
$ws.Cells.Item(138, 10).Value = 30380
$ws.Cells.Item(138, 12).Value = 30380
$ws.Cells.Item(138, 14).Value = -40660
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2273.1304
$ws.Cells.Item(31, 9).Value = 2622.25
$ws.Cells.Item(31, 11).Value = 2622.25
$ws.Cells.Item(31, 13).Value = -2327.25
$ws.Cells.Item(34, 8).Value = 2273.1304
$ws.Cells.Item(34, 9).Value = 2622.25
$ws.Cells.Item(34, 11).Value = 2622.25
$ws.Cells.Item(34, 13).Value = -2420.25
$ws.Cells.Item(134, 8).Value = 20834570
$ws.Cells.Item(134, 9).Value = 1227.5714
$ws.Cells.Item(134, 10).Value = 50001252
$ws.Cells.Item(134, 11).Value = 3682.7142
$ws.Cells.Item(134, 12).Value = 150003756
$ws.Cells.Item(134, 13).Value = -1147.7142
$ws.Cells.Item(134, 14).Value = -150008826
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2554.3845
$ws.Cells.Item(5, 10).Value = 1291.6666
$ws.Cells.Item(5, 12).Value = 3874.9998
$ws.Cells.Item(5, 14).Value = -4098.9998
$ws.Cells.Item(113, 8).Value = 707.4375
$ws.Cells.Item(113, 9).Value = 546.5
$ws.Cells.Item(113, 10).Value = 718.1667
$ws.Cells.Item(113, 11).Value = 1639.5
$ws.Cells.Item(113, 12).Value = 2154.5001
$ws.Cells.Item(113, 13).Value = 530.5
$ws.Cells.Item(113, 14).Value = -6494.5001
$ws.Cells.Item(132, 8).Value = 1876.625
$ws.Cells.Item(132, 9).Value = 1901.6
$ws.Cells.Item(132, 10).Value = 1835
$ws.Cells.Item(132, 11).Value = 17114.4
$ws.Cells.Item(132, 12).Value = 16515
$ws.Cells.Item(132, 13).Value = -14584.4
$ws.Cells.Item(132, 14).Value = -21575
$ws.Cells.Item(134, 8).Value = 3308.8215
$ws.Cells.Item(134, 9).Value = 1429.6428
$ws.Cells.Item(134, 10).Value = 5188
$ws.Cells.Item(134, 11).Value = 4288.928400000001
$ws.Cells.Item(134, 12).Value = 15564
$ws.Cells.Item(134, 13).Value = 781.0715999999993
$ws.Cells.Item(134, 14).Value = -25704
$ws.Cells.Item(135, 8).Value = 2554.3845
$ws.Cells.Item(135, 10).Value = 1291.6666
$ws.Cells.Item(135, 12).Value = 11624.9994
$ws.Cells.Item(135, 14).Value = -16694.9994
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2482.8823
$ws.Cells.Item(80, 9).Value = 1679
$ws.Cells.Item(80, 10).Value = 2921.3635
$ws.Cells.Item(80, 11).Value = 1679
$ws.Cells.Item(80, 12).Value = 2921.3635
$ws.Cells.Item(80, 13).Value = -681
$ws.Cells.Item(80, 14).Value = -4917.363499999999
$ws.Cells.Item(83, 8).Value = 2482.8823
$ws.Cells.Item(83, 9).Value = 1679
$ws.Cells.Item(83, 10).Value = 2921.3635
$ws.Cells.Item(83, 11).Value = 8395
$ws.Cells.Item(83, 12).Value = 14606.8175
$ws.Cells.Item(83, 13).Value = -3403
$ws.Cells.Item(83, 14).Value = -24590.8175
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1821.1072
$ws.Cells.Item(82, 9).Value = 1738.3077
$ws.Cells.Item(82, 10).Value = 1892.8667
$ws.Cells.Item(82, 11).Value = 1738.3077
$ws.Cells.Item(82, 12).Value = 1892.8667
$ws.Cells.Item(82, 13).Value = -1377.3077
$ws.Cells.Item(82, 14).Value = -2614.8667
$ws.Cells.Item(85, 8).Value = 1821.1072
$ws.Cells.Item(85, 9).Value = 1738.3077
$ws.Cells.Item(85, 10).Value = 1892.8667
$ws.Cells.Item(85, 11).Value = 1738.3077
$ws.Cells.Item(85, 12).Value = 1892.8667
$ws.Cells.Item(85, 13).Value = -490.3077000000001
$ws.Cells.Item(85, 14).Value = -4388.8667
$ws.Cells.Item(132, 8).Value = 17034.281
$ws.Cells.Item(132, 9).Value = 1086.5385
$ws.Cells.Item(132, 10).Value = 41912.76
$ws.Cells.Item(132, 11).Value = 3259.6155
$ws.Cells.Item(132, 12).Value = 125738.28
$ws.Cells.Item(132, 13).Value = -729.6155000000003
$ws.Cells.Item(132, 14).Value = -130798.28
$ws.Cells.Item(136, 8).Value = 1053
$ws.Cells.Item(136, 9).Value = 982.58826
$ws.Cells.Item(136, 10).Value = 2250
$ws.Cells.Item(136, 11).Value = 2947.76478
$ws.Cells.Item(136, 12).Value = 6750
$ws.Cells.Item(136, 13).Value = -397.76478
$ws.Cells.Item(136, 14).Value = -11850
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 33340286
$ws.Cells.Item(62, 9).Value = 38468908
$ws.Cells.Item(62, 10).Value = 4250
$ws.Cells.Item(62, 11).Value = 38468908
$ws.Cells.Item(62, 12).Value = 4250
$ws.Cells.Item(62, 13).Value = -38468284
$ws.Cells.Item(62, 14).Value = -5498
$ws.Cells.Item(65, 8).Value = 33340286
$ws.Cells.Item(65, 9).Value = 38468908
$ws.Cells.Item(65, 10).Value = 4250
$ws.Cells.Item(65, 11).Value = 192344540
$ws.Cells.Item(65, 12).Value = 21250
$ws.Cells.Item(65, 13).Value = -192341420
$ws.Cells.Item(65, 14).Value = -27490
$ws.Cells.Item(132, 8).Value = 5911.5
$ws.Cells.Item(132, 9).Value = 6480.0835
$ws.Cells.Item(132, 11).Value = 19440.2505
$ws.Cells.Item(132, 13).Value = -16910.2505
